$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.4458276666666667
$ws.Range("H2").Value = 1.337483
$ws.Range("I2").Value = 0.07239890305200847
$ws.Range("J2").Value = 0.07239890305200847
$ws.Range("M2").Value = 5.197052666666667
$ws.Range("N2").Value = 15.591158
$ws.Range("O2").Value = 0.08925807602505774
$ws.Range("P2").Value = 0.08925807602505774
$ws.Range("Q2").Value = 2.316989863923778
$ws.Range("R2").Value = 20.852908775314
$ws.Range("S2").Value = 0.006462186792746957
$ws.Range("T2").Value = 0.006462186792746957

$ws.Range("G3").Value = 0.4458276666666667
$ws.Range("H3").Value = 1.337483
$ws.Range("I3").Value = 0.07239890305200847
$ws.Range("J3").Value = 0.07239890305200847
$ws.Range("O3").Value = 0.6720990360156781
$ws.Range("P3").Value = 0.6720990360156781
$ws.Range("Q3").Value = 17.44656308258422
$ws.Range("R3").Value = 157.019067743258
$ws.Range("S3").Value = 0.04865923294984743
$ws.Range("T3").Value = 0.04865923294984743

$ws.Range("G4").Value = 0.4458276666666667
$ws.Range("H4").Value = 1.337483
$ws.Range("I4").Value = 0.07239890305200847
$ws.Range("J4").Value = 0.07239890305200847
$ws.Range("O4").Value = 0.2386428879592642
$ws.Range("P4").Value = 0.2386428879592642
$ws.Range("Q4").Value = 6.194768889527555
$ws.Range("R4").Value = 55.752920005748
$ws.Range("S4").Value = 0.01727748330941409
$ws.Range("T4").Value = 0.01727748330941409

$ws.Range("I5").Value = 0.5136151183172798
$ws.Range("J5").Value = 0.5136151183172797
$ws.Range("M5").Value = 5.197052666666667
$ws.Range("N5").Value = 15.591158
$ws.Range("O5").Value = 0.08925807602505774
$ws.Range("P5").Value = 0.08925807602505774
$ws.Range("Q5").Value = 16.43727975055467
$ws.Range("R5").Value = 147.935517754992
$ws.Range("S5").Value = 0.04584429727838279
$ws.Range("T5").Value = 0.04584429727838277

$ws.Range("I6").Value = 0.5136151183172798
$ws.Range("J6").Value = 0.5136151183172797
$ws.Range("O6").Value = 0.6720990360156781
$ws.Range("P6").Value = 0.6720990360156781
$ws.Range("S6").Value = 0.3452002259041222
$ws.Range("T6").Value = 0.3452002259041221

$ws.Range("I7").Value = 0.5136151183172798
$ws.Range("J7").Value = 0.5136151183172797
$ws.Range("O7").Value = 0.2386428879592642
$ws.Range("P7").Value = 0.2386428879592642
$ws.Range("S7").Value = 0.1225705951347748
$ws.Range("T7").Value = 0.1225705951347748

$ws.Range("I8").Value = 0.4139859786307118
$ws.Range("J8").Value = 0.4139859786307118
$ws.Range("M8").Value = 5.197052666666667
$ws.Range("N8").Value = 15.591158
$ws.Range("O8").Value = 0.08925807602505774
$ws.Range("P8").Value = 0.08925807602505774
$ws.Range("Q8").Value = 13.24883770137889
$ws.Range("R8").Value = 119.23953931241
$ws.Range("S8").Value = 0.036951591953928
$ws.Range("T8").Value = 0.036951591953928

$ws.Range("I9").Value = 0.4139859786307118
$ws.Range("J9").Value = 0.4139859786307118
$ws.Range("O9").Value = 0.6720990360156781
$ws.Range("P9").Value = 0.6720990360156781
$ws.Range("Q9").Value = 99.76162879564112
$ws.Range("S9").Value = 0.2782395771617085
$ws.Range("T9").Value = 0.2782395771617085

$ws.Range("I10").Value = 0.4139859786307118
$ws.Range("J10").Value = 0.4139859786307118
$ws.Range("O10").Value = 0.2386428879592642
$ws.Range("P10").Value = 0.2386428879592642
$ws.Range("S10").Value = 0.09879480951507531
$ws.Range("T10").Value = 0.09879480951507531

